$d = $word.ActiveDocument

$replacements = @(
    @{old="144×9=1296"; new="900×2=1800"},
    @{old="114×2=228"; new="499×7=3493"},
    @{old="726×7=5082"; new="354×5=1770"},
    @{old="261×4=1044"; new="564×8=4512"},
    @{old="304×4=1216"; new="177×9=1593"},
    @{old="523×9=4707"; new="956×6=5736"},
    @{old="572×4=2288"; new="890×9=8010"},
    @{old="749×9=6741"; new="373×8=2984"},
    @{old="103×7=721"; new="331×6=1986"},
    @{old="719×9=6471"; new="922×3=2766"},
    @{old="655×4=2620"; new="292×2=584"},
    @{old="521×5=2605"; new="742×4=2968"},
    @{old="237×3=711"; new="526×9=4734"},
    @{old="124×2=248"; new="955×3=2865"},
    @{old="910×5=4550"; new="799×6=4794"},
    @{old="918×8=7344"; new="107×8=856"},
    @{old="194×5=970"; new="477×9=4293"},
    @{old="303×4=1212"; new="343×8=2744"},
    @{old="715×9=6435"; new="958×9=8622"},
    @{old="367×4=1468"; new="266×6=1596"},
    @{old="382×6=2292"; new="112×6=672"},
    @{old="927×5=4635"; new="783×4=3132"},
    @{old="595×7=4165"; new="589×3=1767"},
    @{old="722×8=5776"; new="809×6=4854"},
    @{old="880×4=3520"; new="718×5=3590"}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.new, 2)
}
